# May2025.xlsx — add a "Categories" column (D) and move the dollar amount
# into a new "Amount" column (E): Categories/Food/Transportation/Clothing/
# Subscriptions/Miscellaneous/Beauty & Care/Amount are new shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy column D's cell format into column E for every data row, and
#        for the header row, so the new column inherits the right number
#        format / borders (same style index as D previously had). ---
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("D4:D37").Copy()
$ws.Range("E4:E37").PasteSpecial(-4122)

# --- 2. Header row: rename D3 to "Categories" (was "$  Amount"); the new
#        "Amount" header text goes into E3 (written later, after the
#        category column, to match the original author's edit order). ---
$ws.Range("D3").Value = "Categories"

# --- 3. Fill in the new Categories column, grouped by category (this
#        mirrors the order the categories were first introduced into the
#        shared-string table: Food, Transportation, Clothing,
#        Subscriptions, Miscellaneous, Beauty & Care). ---
$foodRows          = @(4,5,11,12,14,16,18,19,20,21,22,25,26,27,28,30,31,33,34,35,37)
$transportRows     = @(6,13)
$clothingRows      = @(24,29,32,36)
$subscriptionRows  = @(10)
$miscRows          = @(7,8,9,15,17)
$beautyRows        = @(23)

foreach ($r in $foodRows)         { $ws.Cells.Item($r, 4).Value = "Food" }
foreach ($r in $transportRows)    { $ws.Cells.Item($r, 4).Value = "Transportation" }
foreach ($r in $clothingRows)     { $ws.Cells.Item($r, 4).Value = "Clothing" }
foreach ($r in $subscriptionRows) { $ws.Cells.Item($r, 4).Value = "Subscriptions" }
foreach ($r in $miscRows)         { $ws.Cells.Item($r, 4).Value = "Miscellaneous" }
foreach ($r in $beautyRows)       { $ws.Cells.Item($r, 4).Value = "Beauty & Care" }

# --- 4. Move each row's dollar amount (previously in D) into the new
#        Amount column E. ---
$amounts = @{
    4  = 2.5
    5  = 12.52
    6  = 7.96
    7  = 2
    8  = 2
    9  = 10
    10 = 20
    11 = 21.85
    12 = 2.5
    13 = 55
    14 = 7.17
    15 = 4.28
    16 = 9.67
    17 = 4.09
    18 = 8.29
    19 = 13.99
    20 = 5.87
    21 = 2.26
    22 = 4.51
    23 = 35.81
    24 = 41.49
    25 = 4.02
    26 = 9.43
    27 = 10.19
    28 = 10.49
    29 = 22.21
    30 = 12.46
    31 = 19.39
    32 = 22.21
    33 = 5.93
    34 = 19.39
    35 = 12.01
    36 = 12.53
    37 = 12.04
}
for ($r = 4; $r -le 37; $r++) {
    $ws.Cells.Item($r, 5).Value = $amounts[$r]
}

# --- 5. "Amount" header label goes in E3, written last so it lands at the
#        end of the shared-string table. ---
$ws.Range("E3").Value = "Amount"

# --- 6. Column E width + sheet view (zoom / selection) cosmetic updates. ---
$ws.Columns.Item(5).ColumnWidth = 21.83203125

$excel.ActiveWindow.Zoom = 138
$ws.Range("D8").Select()

Write-Output "done"
